# Edit script for unitTest_repeatUntil.xlsx
# Adds new rdbms/localdb/web commands and replaces the tn.5250-profile
# command set on column Z with the new step.inTime command set,
# per the commit: new assertResultMatch/assertResultNotMatch rdbms commands,
# plus supporting new commands across several system command lists.

$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

# ---- column helpers (1-based column numbers) ----
# A=1 O=15 T=20 Z=26 AA=27

# ---- 1) rdbms command list -> column T, rows 2-9 ----
$rdbms = @(
    'assertResultMatch(var,columns,search)',
    'assertResultNotMatch(var,columns,search)',
    'resultToCSV(var,csvFile,delim,showHeader)',
    'runFile(var,db,file)',
    'runSQL(var,db,sql)',
    'runSQLs(var,db,sqls)',
    'saveResult(db,sql,output)',
    'saveResults(db,sqls,outputDir)'
)
for ($i = 0; $i -lt $rdbms.Length; $i++) {
    $sys.Cells.Item(2 + $i, 20).Value = $rdbms[$i]
}

# ---- 2) localdb command list -> column O, rows 2-13 ----
$localdb = @(
    'cloneTable(var,source,target)',
    'dropTables(var,tables)',
    'exportCSV(sql,output)',
    'exportEXCEL(sql,output,sheet,start)',
    'exportJSON(sql,output,header)',
    'exportXML(sql,output,root,row,cell)',
    'importCSV(var,csv,table)',
    'importEXCEL(var,excel,sheet,ranges,table)',
    'importRecords(var,sourceDb,sql,table)',
    'purge(var)',
    'queryAsCSV(var,sql)',
    'runSQLs(var,sqls)'
)
for ($i = 0; $i -lt $localdb.Length; $i++) {
    $sys.Cells.Item(2 + $i, 15).Value = $localdb[$i]
}

# ---- 3) step.inTime command list -> column Z, rows 2-4 (replaces the old
#        tn.5250 profile commands; rows 5-6 are cleared) ----
$sys.Cells.Item(1, 26).Value = "step.inTime"
$stepInTime = @(
    'observe(prompt,waitMs)',
    'perform(instructions,waitMs)',
    'validate(prompt,responses,passResponses,waitMs)'
)
for ($i = 0; $i -lt $stepInTime.Length; $i++) {
    $sys.Cells.Item(2 + $i, 26).Value = $stepInTime[$i]
}
$sys.Cells.Item(5, 26).ClearContents()
$sys.Cells.Item(6, 26).ClearContents()

# ---- 4) "target" category list -> column A, row 26 relabeled ----
$sys.Cells.Item(26, 1).Value = "step.inTime"

# ---- 5) web command list -> column AA, rows 2-145 ----
$web = @(
    'assertAndClick(locator,label)',
    'assertAttribute(locator,attrName,value)',
    'assertAttributeContain(locator,attrName,contains)',
    'assertAttributeNotContain(locator,attrName,contains)',
    'assertAttributeNotPresent(locator,attrName)',
    'assertAttributePresent(locator,attrName)',
    'assertChecked(locator)',
    'assertContainCount(locator,text,count)',
    'assertCssNotPresent(locator,property)',
    'assertCssPresent(locator,property,value)',
    'assertElementByAttributes(nameValues)',
    'assertElementByText(locator,text)',
    'assertElementCount(locator,count)',
    'assertElementDisabled(locator)',
    'assertElementEnabled(locator)',
    'assertElementNotPresent(locator)',
    'assertElementPresent(locator)',
    'assertElementsPresent(prefix)',
    'assertFocus(locator)',
    'assertFrameCount(count)',
    'assertFramePresent(frameName)',
    'assertIECompatMode()',
    'assertIENativeMode()',
    'assertLinkByLabel(label)',
    'assertMultiSelect(locator)',
    'assertNotChecked(locator)',
    'assertNotFocus(locator)',
    'assertNotText(locator,text)',
    'assertNotVisible(locator)',
    'assertOneMatch(locator)',
    'assertScrollbarHNotPresent(locator)',
    'assertScrollbarHPresent(locator)',
    'assertScrollbarVNotPresent(locator)',
    'assertScrollbarVPresent(locator)',
    'assertSingleSelect(locator)',
    'assertTable(locator,row,column,text)',
    'assertText(locator,text)',
    'assertTextContains(locator,text)',
    'assertTextCount(locator,text,count)',
    'assertTextList(locator,list,ignoreOrder)',
    'assertTextMatches(text,minMatch,scrollTo)',
    'assertTextNotContain(locator,text)',
    'assertTextNotPresent(text)',
    'assertTextOrder(locator,descending)',
    'assertTextPresent(text)',
    'assertTitle(text)',
    'assertValue(locator,value)',
    'assertValueOrder(locator,descending)',
    'assertVisible(locator)',
    'checkAll(locator,waitMs)',
    'clearLocalStorage()',
    'click(locator)',
    'clickAll(locator)',
    'clickAndWait(locator,waitMs)',
    'clickByLabel(label)',
    'clickByLabelAndWait(label,waitMs)',
    'clickOffset(locator,x,y)',
    'clickWithKeys(locator,keys)',
    'close()',
    'closeAll()',
    'deselect(locator,text)',
    'deselectMulti(locator,array)',
    'dismissInvalidCert()',
    'dismissInvalidCertPopup()',
    'doubleClick(locator)',
    'doubleClickAndWait(locator,waitMs)',
    'doubleClickByLabel(label)',
    'doubleClickByLabelAndWait(label,waitMs)',
    'dragAndDrop(fromLocator,toLocator)',
    'dragTo(fromLocator,xOffset,yOffset)',
    'editLocalStorage(key,value)',
    'executeScript(var,script)',
    'focus(locator)',
    'goBack()',
    'goBackAndWait()',
    'maximizeWindow()',
    'mouseOver(locator)',
    'open(url)',
    'openAndWait(url,waitMs)',
    'openHttpBasic(url,username,password)',
    'openIgnoreTimeout(url)',
    'refresh()',
    'refreshAndWait()',
    'resizeWindow(width,height)',
    'rightClick(locator)',
    'saveAllWindowIds(var)',
    'saveAllWindowNames(var)',
    'saveAttribute(var,locator,attrName)',
    'saveAttributeList(var,locator,attrName)',
    'saveBrowserVersion(var)',
    'saveCount(var,locator)',
    'saveDivsAsCsv(headers,rows,cells,nextPage,file)',
    'saveElement(var,locator)',
    'saveElements(var,locator)',
    'saveInfiniteDivsAsCsv(config,file)',
    'saveInfiniteTableAsCsv(config,file)',
    'saveLocalStorage(var,key)',
    'saveLocation(var)',
    'savePageAs(var,sessionIdName,url)',
    'savePageAsFile(sessionIdName,url,file)',
    'saveSelectedText(var,locator)',
    'saveSelectedValue(var,locator)',
    'saveTableAsCsv(locator,nextPageLocator,file)',
    'saveText(var,locator)',
    'saveTextArray(var,locator)',
    'saveTextSubstringAfter(var,locator,delim)',
    'saveTextSubstringBefore(var,locator,delim)',
    'saveTextSubstringBetween(var,locator,start,end)',
    'saveTitle(var)',
    'saveValue(var,locator)',
    'saveValues(var,locator)',
    'screenshot(file,locator)',
    'scrollElement(locator,xOffset,yOffset)',
    'scrollLeft(locator,pixel)',
    'scrollPage(xOffset,yOffset)',
    'scrollRight(locator,pixel)',
    'scrollTo(locator)',
    'select(locator,text)',
    'selectAllOptions(locator)',
    'selectFrame(locator)',
    'selectMulti(locator,array)',
    'selectMultiByValue(locator,array)',
    'selectMultiOptions(locator)',
    'selectText(locator)',
    'selectWindow(winId)',
    'selectWindowAndWait(winId,waitMs)',
    'selectWindowByIndex(index)',
    'selectWindowByIndexAndWait(index,waitMs)',
    'switchBrowser(profile,config)',
    'toggleSelections(locator)',
    'type(locator,value)',
    'typeKeys(locator,value)',
    'uncheckAll(locator,waitMs)',
    'unselectAllText()',
    'updateAttribute(locator,attrName,value)',
    'upload(fieldLocator,file)',
    'verifyContainText(locator,text)',
    'verifyText(locator,text)',
    'wait(waitMs)',
    'waitForElementPresent(locator)',
    'waitForElementsPresent(locators)',
    'waitForPopUp(winId,waitMs)',
    'waitForTextPresent(text)',
    'waitForTitle(text)'
)
for ($i = 0; $i -lt $web.Length; $i++) {
    $sys.Cells.Item(2 + $i, 27).Value = $web[$i]
}

# ---- 6) defined names: widen existing ranges, add the new one ----
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("localdb").RefersTo = "='#system'!`$O`$2:`$O`$13"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$145"
$wb.Names.Add("step.inTime", "='#system'!`$Z`$2:`$Z`$4")

Write-Host "Edit complete."
